$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 2.401444666666666
$ws.Range("H2").Value = 7.204333999999999
$ws.Range("I2").Value = 0.5723125574599716
$ws.Range("J2").Value = 0.5723125574599716
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 3.778439
$ws.Range("N2").Value = 11.335317
$ws.Range("O2").Value = 0.4252971528324392
$ws.Range("P2").Value = 0.4252971528324392
$ws.Range("Q2").Value = 9.073712184875333
$ws.Range("R2").Value = 81.66340966387799
$ws.Range("S2").Value = 0.2434029012179777
$ws.Range("T2").Value = 0.2434029012179777

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 2.401444666666666
$ws.Range("H3").Value = 7.204333999999999
$ws.Range("I3").Value = 0.5723125574599716
$ws.Range("J3").Value = 0.5723125574599716
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 4.333403333333333
$ws.Range("N3").Value = 13.00021
$ws.Range("O3").Value = 0.4877633593505858
$ws.Range("P3").Value = 0.4877633593505858
$ws.Range("Q3").Value = 10.40642832334889
$ws.Range("R3").Value = 93.65785491013999
$ws.Range("S3").Value = 0.2791530956252009
$ws.Range("T3").Value = 0.2791530956252009

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 2.401444666666666
$ws.Range("H4").Value = 7.204333999999999
$ws.Range("I4").Value = 0.5723125574599716
$ws.Range("J4").Value = 0.5723125574599716
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 0.2909853333333334
$ws.Range("N4").Value = 0.8729560000000001
$ws.Range("O4").Value = 0.03275300561492853
$ws.Range("P4").Value = 0.03275300561492853
$ws.Range("Q4").Value = 0.6987851768115556
$ws.Range("R4").Value = 6.289066591304
$ws.Range("S4").Value = 0.01874495640798056
$ws.Range("T4").Value = 0.01874495640798056

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 2.401444666666666
$ws.Range("H5").Value = 7.204333999999999
$ws.Range("I5").Value = 0.5723125574599716
$ws.Range("J5").Value = 0.5723125574599716
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 0.4814053333333333
$ws.Range("N5").Value = 1.444216
$ws.Range("O5").Value = 0.0541864822020464
$ws.Range("P5").Value = 0.05418648220204641
$ws.Range("Q5").Value = 1.156068270238222
$ws.Range("R5").Value = 10.404614432144
$ws.Range("S5").Value = 0.03101160420881241
$ws.Range("T5").Value = 0.03101160420881241

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 0.9802360000000001
$ws.Range("H6").Value = 2.940708
$ws.Range("I6").Value = 0.2336099514851752
$ws.Range("J6").Value = 0.2336099514851752
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 3.778439
$ws.Range("N6").Value = 11.335317
$ws.Range("O6").Value = 0.4252971528324392
$ws.Range("P6").Value = 0.4252971528324392
$ws.Range("Q6").Value = 3.703761931604001
$ws.Range("R6").Value = 33.333857384436
$ws.Range("S6").Value = 0.09935364723996928
$ws.Range("T6").Value = 0.09935364723996928

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 0.9802360000000001
$ws.Range("H7").Value = 2.940708
$ws.Range("I7").Value = 0.2336099514851752
$ws.Range("J7").Value = 0.2336099514851752
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 4.333403333333333
$ws.Range("N7").Value = 13.00021
$ws.Range("O7").Value = 0.4877633593505858
$ws.Range("P7").Value = 0.4877633593505858
$ws.Range("Q7").Value = 4.247757949853334
$ws.Range("R7").Value = 38.22982154868
$ws.Range("S7").Value = 0.1139463747141364
$ws.Range("T7").Value = 0.1139463747141364

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 0.9802360000000001
$ws.Range("H8").Value = 2.940708
$ws.Range("I8").Value = 0.2336099514851752
$ws.Range("J8").Value = 0.2336099514851752
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 0.2909853333333334
$ws.Range("N8").Value = 0.8729560000000001
$ws.Range("O8").Value = 0.03275300561492853
$ws.Range("P8").Value = 0.03275300561492853
$ws.Range("Q8").Value = 0.2852342992053334
$ws.Range("R8").Value = 2.567108692848
$ws.Range("S8").Value = 0.007651428052697126
$ws.Range("T8").Value = 0.007651428052697126

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 0.9802360000000001
$ws.Range("H9").Value = 2.940708
$ws.Range("I9").Value = 0.2336099514851752
$ws.Range("J9").Value = 0.2336099514851752
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 0.4814053333333333
$ws.Range("N9").Value = 1.444216
$ws.Range("O9").Value = 0.0541864822020464
$ws.Range("P9").Value = 0.05418648220204641
$ws.Range("Q9").Value = 0.4718908383253334
$ws.Range("R9").Value = 4.247017544928
$ws.Range("S9").Value = 0.01265850147837237
$ws.Range("T9").Value = 0.01265850147837237

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 0.8143563333333333
$ws.Range("H10").Value = 2.443069
$ws.Range("I10").Value = 0.1940774910548533
$ws.Range("J10").Value = 0.1940774910548533
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 3.778439
$ws.Range("N10").Value = 11.335317
$ws.Range("O10").Value = 0.4252971528324392
$ws.Range("P10").Value = 0.4252971528324392
$ws.Range("Q10").Value = 3.076995729763667
$ws.Range("R10").Value = 27.692961567873
$ws.Range("S10").Value = 0.08254060437449229
$ws.Range("T10").Value = 0.08254060437449229

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 0.8143563333333333
$ws.Range("H11").Value = 2.443069
$ws.Range("I11").Value = 0.1940774910548533
$ws.Range("J11").Value = 0.1940774910548533
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 4.333403333333333
$ws.Range("N11").Value = 13.00021
$ws.Range("O11").Value = 0.4877633593505858
$ws.Range("P11").Value = 0.4877633593505858
$ws.Range("Q11").Value = 3.528934449387778
$ws.Range("R11").Value = 31.76041004449
$ws.Range("S11").Value = 0.0946638890112485
$ws.Range("T11").Value = 0.09466388901124852

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 0.8143563333333333
$ws.Range("H12").Value = 2.443069
$ws.Range("I12").Value = 0.1940774910548533
$ws.Range("J12").Value = 0.1940774910548533
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 0.2909853333333334
$ws.Range("N12").Value = 0.8729560000000001
$ws.Range("O12").Value = 0.03275300561492853
$ws.Range("P12").Value = 0.03275300561492853
$ws.Range("Q12").Value = 0.2369657491071112
$ws.Range("R12").Value = 2.132691741964
$ws.Range("S12").Value = 0.006356621154250851
$ws.Range("T12").Value = 0.006356621154250851

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 0.8143563333333333
$ws.Range("H13").Value = 2.443069
$ws.Range("I13").Value = 0.1940774910548533
$ws.Range("J13").Value = 0.1940774910548533
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 0.4814053333333333
$ws.Range("N13").Value = 1.444216
$ws.Range("O13").Value = 0.0541864822020464
$ws.Range("P13").Value = 0.05418648220204641
$ws.Range("Q13").Value = 0.3920354821004444
$ws.Range("R13").Value = 3.528319338904
$ws.Range("S13").Value = 0.01051637651486163
$ws.Range("T13").Value = 0.01051637651486163

